# Append new sensor-log rows generated on 2026-01-28 to the PIR, Humidity,
# and Temperature sheets, matching the automated logger export.

$wb = $excel.ActiveWorkbook

$wsPIR = $wb.Worksheets.Item("PIR")
$wsHumidity = $wb.Worksheets.Item("Humidity")
$wsTemperature = $wb.Worksheets.Item("Temperature")

# --- PIR sheet: append rows 33-45 ---
$wsPIR.Range("A33").NumberFormat = "@"
$wsPIR.Range("A33").Value = "2026-01-28"
$wsPIR.Range("B33").Value = "17:04:37"
$wsPIR.Range("C33").Value = "17:00"
$wsPIR.Range("D33").Value = "Bathroom"
$wsPIR.Range("E33").Value = "No Motion"
$wsPIR.Range("F33").Value = "Inactive"

$wsPIR.Range("A34").NumberFormat = "@"
$wsPIR.Range("A34").Value = "2026-01-28"
$wsPIR.Range("B34").Value = "17:04:38"
$wsPIR.Range("C34").Value = "17:00"
$wsPIR.Range("D34").Value = "Bathroom"
$wsPIR.Range("E34").Value = "No Motion"
$wsPIR.Range("F34").Value = "Inactive"

$wsPIR.Range("A35").NumberFormat = "@"
$wsPIR.Range("A35").Value = "2026-01-28"
$wsPIR.Range("B35").Value = "17:04:43"
$wsPIR.Range("C35").Value = "17:00"
$wsPIR.Range("D35").Value = "Bathroom"
$wsPIR.Range("E35").Value = "No Motion"
$wsPIR.Range("F35").Value = "Inactive"

$wsPIR.Range("A36").NumberFormat = "@"
$wsPIR.Range("A36").Value = "2026-01-28"
$wsPIR.Range("B36").Value = "17:04:48"
$wsPIR.Range("C36").Value = "17:00"
$wsPIR.Range("D36").Value = "Bathroom"
$wsPIR.Range("E36").Value = "No Motion"
$wsPIR.Range("F36").Value = "Inactive"

$wsPIR.Range("A37").NumberFormat = "@"
$wsPIR.Range("A37").Value = "2026-01-28"
$wsPIR.Range("B37").Value = "17:04:53"
$wsPIR.Range("C37").Value = "17:00"
$wsPIR.Range("D37").Value = "Bathroom"
$wsPIR.Range("E37").Value = "No Motion"
$wsPIR.Range("F37").Value = "Inactive"

$wsPIR.Range("A38").NumberFormat = "@"
$wsPIR.Range("A38").Value = "2026-01-28"
$wsPIR.Range("B38").Value = "17:04:58"
$wsPIR.Range("C38").Value = "17:00"
$wsPIR.Range("D38").Value = "Bathroom"
$wsPIR.Range("E38").Value = "No Motion"
$wsPIR.Range("F38").Value = "Inactive"

$wsPIR.Range("A39").NumberFormat = "@"
$wsPIR.Range("A39").Value = "2026-01-28"
$wsPIR.Range("B39").Value = "17:05:03"
$wsPIR.Range("C39").Value = "17:00"
$wsPIR.Range("D39").Value = "Bathroom"
$wsPIR.Range("E39").Value = "No Motion"
$wsPIR.Range("F39").Value = "Inactive"

$wsPIR.Range("A40").NumberFormat = "@"
$wsPIR.Range("A40").Value = "2026-01-28"
$wsPIR.Range("B40").Value = "17:05:08"
$wsPIR.Range("C40").Value = "17:00"
$wsPIR.Range("D40").Value = "Bathroom"
$wsPIR.Range("E40").Value = "No Motion"
$wsPIR.Range("F40").Value = "Inactive"

$wsPIR.Range("A41").NumberFormat = "@"
$wsPIR.Range("A41").Value = "2026-01-28"
$wsPIR.Range("B41").Value = "17:05:13"
$wsPIR.Range("C41").Value = "17:00"
$wsPIR.Range("D41").Value = "Bathroom"
$wsPIR.Range("E41").Value = "No Motion"
$wsPIR.Range("F41").Value = "Inactive"

$wsPIR.Range("A42").NumberFormat = "@"
$wsPIR.Range("A42").Value = "2026-01-28"
$wsPIR.Range("B42").Value = "17:05:18"
$wsPIR.Range("C42").Value = "17:00"
$wsPIR.Range("D42").Value = "Bathroom"
$wsPIR.Range("E42").Value = "No Motion"
$wsPIR.Range("F42").Value = "Inactive"

$wsPIR.Range("A43").NumberFormat = "@"
$wsPIR.Range("A43").Value = "2026-01-28"
$wsPIR.Range("B43").Value = "17:05:23"
$wsPIR.Range("C43").Value = "17:00"
$wsPIR.Range("D43").Value = "Bathroom"
$wsPIR.Range("E43").Value = "No Motion"
$wsPIR.Range("F43").Value = "Inactive"

$wsPIR.Range("A44").NumberFormat = "@"
$wsPIR.Range("A44").Value = "2026-01-28"
$wsPIR.Range("B44").Value = "17:05:28"
$wsPIR.Range("C44").Value = "17:00"
$wsPIR.Range("D44").Value = "Bathroom"
$wsPIR.Range("E44").Value = "No Motion"
$wsPIR.Range("F44").Value = "Inactive"

$wsPIR.Range("A45").NumberFormat = "@"
$wsPIR.Range("A45").Value = "2026-01-28"
$wsPIR.Range("B45").Value = "17:05:33"
$wsPIR.Range("C45").Value = "17:00"
$wsPIR.Range("D45").Value = "Bathroom"
$wsPIR.Range("E45").Value = "No Motion"
$wsPIR.Range("F45").Value = "Inactive"

# --- Humidity sheet: append rows 32-45 ---
$wsHumidity.Range("A32").NumberFormat = "@"
$wsHumidity.Range("E32").NumberFormat = "@"
$wsHumidity.Range("A32").Value = "2026-01-28"
$wsHumidity.Range("B32").Value = "17:04:37"
$wsHumidity.Range("C32").Value = "17:00"
$wsHumidity.Range("D32").Value = "Bathroom"
$wsHumidity.Range("E32").Value = "86.6%"
$wsHumidity.Range("F32").Value = "Active"

$wsHumidity.Range("A33").NumberFormat = "@"
$wsHumidity.Range("E33").NumberFormat = "@"
$wsHumidity.Range("A33").Value = "2026-01-28"
$wsHumidity.Range("B33").Value = "17:04:37"
$wsHumidity.Range("C33").Value = "17:00"
$wsHumidity.Range("D33").Value = "Bathroom"
$wsHumidity.Range("E33").Value = "87.5%"
$wsHumidity.Range("F33").Value = "Active"

$wsHumidity.Range("A34").NumberFormat = "@"
$wsHumidity.Range("E34").NumberFormat = "@"
$wsHumidity.Range("A34").Value = "2026-01-28"
$wsHumidity.Range("B34").Value = "17:04:39"
$wsHumidity.Range("C34").Value = "17:00"
$wsHumidity.Range("D34").Value = "Bathroom"
$wsHumidity.Range("E34").Value = "86.6%"
$wsHumidity.Range("F34").Value = "Active"

$wsHumidity.Range("A35").NumberFormat = "@"
$wsHumidity.Range("E35").NumberFormat = "@"
$wsHumidity.Range("A35").Value = "2026-01-28"
$wsHumidity.Range("B35").Value = "17:04:43"
$wsHumidity.Range("C35").Value = "17:00"
$wsHumidity.Range("D35").Value = "Bathroom"
$wsHumidity.Range("E35").Value = "87.5%"
$wsHumidity.Range("F35").Value = "Active"

$wsHumidity.Range("A36").NumberFormat = "@"
$wsHumidity.Range("E36").NumberFormat = "@"
$wsHumidity.Range("A36").Value = "2026-01-28"
$wsHumidity.Range("B36").Value = "17:04:47"
$wsHumidity.Range("C36").Value = "17:00"
$wsHumidity.Range("D36").Value = "Bathroom"
$wsHumidity.Range("E36").Value = "87.5%"
$wsHumidity.Range("F36").Value = "Active"

$wsHumidity.Range("A37").NumberFormat = "@"
$wsHumidity.Range("E37").NumberFormat = "@"
$wsHumidity.Range("A37").Value = "2026-01-28"
$wsHumidity.Range("B37").Value = "17:04:51"
$wsHumidity.Range("C37").Value = "17:00"
$wsHumidity.Range("D37").Value = "Bathroom"
$wsHumidity.Range("E37").Value = "86.6%"
$wsHumidity.Range("F37").Value = "Active"

$wsHumidity.Range("A38").NumberFormat = "@"
$wsHumidity.Range("E38").NumberFormat = "@"
$wsHumidity.Range("A38").Value = "2026-01-28"
$wsHumidity.Range("B38").Value = "17:04:55"
$wsHumidity.Range("C38").Value = "17:00"
$wsHumidity.Range("D38").Value = "Bathroom"
$wsHumidity.Range("E38").Value = "87.6%"
$wsHumidity.Range("F38").Value = "Active"

$wsHumidity.Range("A39").NumberFormat = "@"
$wsHumidity.Range("E39").NumberFormat = "@"
$wsHumidity.Range("A39").Value = "2026-01-28"
$wsHumidity.Range("B39").Value = "17:05:07"
$wsHumidity.Range("C39").Value = "17:00"
$wsHumidity.Range("D39").Value = "Bathroom"
$wsHumidity.Range("E39").Value = "87.6%"
$wsHumidity.Range("F39").Value = "Active"

$wsHumidity.Range("A40").NumberFormat = "@"
$wsHumidity.Range("E40").NumberFormat = "@"
$wsHumidity.Range("A40").Value = "2026-01-28"
$wsHumidity.Range("B40").Value = "17:05:11"
$wsHumidity.Range("C40").Value = "17:00"
$wsHumidity.Range("D40").Value = "Bathroom"
$wsHumidity.Range("E40").Value = "86.7%"
$wsHumidity.Range("F40").Value = "Active"

$wsHumidity.Range("A41").NumberFormat = "@"
$wsHumidity.Range("E41").NumberFormat = "@"
$wsHumidity.Range("A41").Value = "2026-01-28"
$wsHumidity.Range("B41").Value = "17:05:15"
$wsHumidity.Range("C41").Value = "17:00"
$wsHumidity.Range("D41").Value = "Bathroom"
$wsHumidity.Range("E41").Value = "87.6%"
$wsHumidity.Range("F41").Value = "Active"

$wsHumidity.Range("A42").NumberFormat = "@"
$wsHumidity.Range("E42").NumberFormat = "@"
$wsHumidity.Range("A42").Value = "2026-01-28"
$wsHumidity.Range("B42").Value = "17:05:19"
$wsHumidity.Range("C42").Value = "17:00"
$wsHumidity.Range("D42").Value = "Bathroom"
$wsHumidity.Range("E42").Value = "87.6%"
$wsHumidity.Range("F42").Value = "Active"

$wsHumidity.Range("A43").NumberFormat = "@"
$wsHumidity.Range("E43").NumberFormat = "@"
$wsHumidity.Range("A43").Value = "2026-01-28"
$wsHumidity.Range("B43").Value = "17:05:23"
$wsHumidity.Range("C43").Value = "17:00"
$wsHumidity.Range("D43").Value = "Bathroom"
$wsHumidity.Range("E43").Value = "87.6%"
$wsHumidity.Range("F43").Value = "Active"

$wsHumidity.Range("A44").NumberFormat = "@"
$wsHumidity.Range("E44").NumberFormat = "@"
$wsHumidity.Range("A44").Value = "2026-01-28"
$wsHumidity.Range("B44").Value = "17:05:27"
$wsHumidity.Range("C44").Value = "17:00"
$wsHumidity.Range("D44").Value = "Bathroom"
$wsHumidity.Range("E44").Value = "87.6%"
$wsHumidity.Range("F44").Value = "Active"

$wsHumidity.Range("A45").NumberFormat = "@"
$wsHumidity.Range("E45").NumberFormat = "@"
$wsHumidity.Range("A45").Value = "2026-01-28"
$wsHumidity.Range("B45").Value = "17:05:35"
$wsHumidity.Range("C45").Value = "17:00"
$wsHumidity.Range("D45").Value = "Bathroom"
$wsHumidity.Range("E45").Value = "87.6%"
$wsHumidity.Range("F45").Value = "Active"

# --- Temperature sheet: append rows 32-45 ---
$wsTemperature.Range("A32").NumberFormat = "@"
$wsTemperature.Range("A32").Value = "2026-01-28"
$wsTemperature.Range("B32").Value = "17:04:37"
$wsTemperature.Range("C32").Value = "17:00"
$wsTemperature.Range("D32").Value = "Bathroom"
$wsTemperature.Range("E32").Value = "22.9C"
$wsTemperature.Range("F32").Value = "Active"

$wsTemperature.Range("A33").NumberFormat = "@"
$wsTemperature.Range("A33").Value = "2026-01-28"
$wsTemperature.Range("B33").Value = "17:04:37"
$wsTemperature.Range("C33").Value = "17:00"
$wsTemperature.Range("D33").Value = "Bathroom"
$wsTemperature.Range("E33").Value = "22.9C"
$wsTemperature.Range("F33").Value = "Active"

$wsTemperature.Range("A34").NumberFormat = "@"
$wsTemperature.Range("A34").Value = "2026-01-28"
$wsTemperature.Range("B34").Value = "17:04:39"
$wsTemperature.Range("C34").Value = "17:00"
$wsTemperature.Range("D34").Value = "Bathroom"
$wsTemperature.Range("E34").Value = "22.8C"
$wsTemperature.Range("F34").Value = "Active"

$wsTemperature.Range("A35").NumberFormat = "@"
$wsTemperature.Range("A35").Value = "2026-01-28"
$wsTemperature.Range("B35").Value = "17:04:43"
$wsTemperature.Range("C35").Value = "17:00"
$wsTemperature.Range("D35").Value = "Bathroom"
$wsTemperature.Range("E35").Value = "22.8C"
$wsTemperature.Range("F35").Value = "Active"

$wsTemperature.Range("A36").NumberFormat = "@"
$wsTemperature.Range("A36").Value = "2026-01-28"
$wsTemperature.Range("B36").Value = "17:04:47"
$wsTemperature.Range("C36").Value = "17:00"
$wsTemperature.Range("D36").Value = "Bathroom"
$wsTemperature.Range("E36").Value = "22.8C"
$wsTemperature.Range("F36").Value = "Active"

$wsTemperature.Range("A37").NumberFormat = "@"
$wsTemperature.Range("A37").Value = "2026-01-28"
$wsTemperature.Range("B37").Value = "17:04:51"
$wsTemperature.Range("C37").Value = "17:00"
$wsTemperature.Range("D37").Value = "Bathroom"
$wsTemperature.Range("E37").Value = "22.8C"
$wsTemperature.Range("F37").Value = "Active"

$wsTemperature.Range("A38").NumberFormat = "@"
$wsTemperature.Range("A38").Value = "2026-01-28"
$wsTemperature.Range("B38").Value = "17:04:55"
$wsTemperature.Range("C38").Value = "17:00"
$wsTemperature.Range("D38").Value = "Bathroom"
$wsTemperature.Range("E38").Value = "22.8C"
$wsTemperature.Range("F38").Value = "Active"

$wsTemperature.Range("A39").NumberFormat = "@"
$wsTemperature.Range("A39").Value = "2026-01-28"
$wsTemperature.Range("B39").Value = "17:05:07"
$wsTemperature.Range("C39").Value = "17:00"
$wsTemperature.Range("D39").Value = "Bathroom"
$wsTemperature.Range("E39").Value = "22.9C"
$wsTemperature.Range("F39").Value = "Active"

$wsTemperature.Range("A40").NumberFormat = "@"
$wsTemperature.Range("A40").Value = "2026-01-28"
$wsTemperature.Range("B40").Value = "17:05:11"
$wsTemperature.Range("C40").Value = "17:00"
$wsTemperature.Range("D40").Value = "Bathroom"
$wsTemperature.Range("E40").Value = "22.9C"
$wsTemperature.Range("F40").Value = "Active"

$wsTemperature.Range("A41").NumberFormat = "@"
$wsTemperature.Range("A41").Value = "2026-01-28"
$wsTemperature.Range("B41").Value = "17:05:15"
$wsTemperature.Range("C41").Value = "17:00"
$wsTemperature.Range("D41").Value = "Bathroom"
$wsTemperature.Range("E41").Value = "22.8C"
$wsTemperature.Range("F41").Value = "Active"

$wsTemperature.Range("A42").NumberFormat = "@"
$wsTemperature.Range("A42").Value = "2026-01-28"
$wsTemperature.Range("B42").Value = "17:05:19"
$wsTemperature.Range("C42").Value = "17:00"
$wsTemperature.Range("D42").Value = "Bathroom"
$wsTemperature.Range("E42").Value = "22.8C"
$wsTemperature.Range("F42").Value = "Active"

$wsTemperature.Range("A43").NumberFormat = "@"
$wsTemperature.Range("A43").Value = "2026-01-28"
$wsTemperature.Range("B43").Value = "17:05:24"
$wsTemperature.Range("C43").Value = "17:00"
$wsTemperature.Range("D43").Value = "Bathroom"
$wsTemperature.Range("E43").Value = "22.9C"
$wsTemperature.Range("F43").Value = "Active"

$wsTemperature.Range("A44").NumberFormat = "@"
$wsTemperature.Range("A44").Value = "2026-01-28"
$wsTemperature.Range("B44").Value = "17:05:28"
$wsTemperature.Range("C44").Value = "17:00"
$wsTemperature.Range("D44").Value = "Bathroom"
$wsTemperature.Range("E44").Value = "22.9C"
$wsTemperature.Range("F44").Value = "Active"

$wsTemperature.Range("A45").NumberFormat = "@"
$wsTemperature.Range("A45").Value = "2026-01-28"
$wsTemperature.Range("B45").Value = "17:05:36"
$wsTemperature.Range("C45").Value = "17:00"
$wsTemperature.Range("D45").Value = "Bathroom"
$wsTemperature.Range("E45").Value = "22.9C"
$wsTemperature.Range("F45").Value = "Active"
